$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.975.42'
$ws.Range('E2').Value = '  +0.35%  '

$ws.Range('D3').Value = '1.594.13'
$ws.Range('E3').Value = '  +0.29%  '

$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('E5').Value = '  +0.26%  '

$ws.Range('E6').Value = '  -0.17%  '

$ws.Range('E7').Value = '  -0.29%  '

$ws.Range('D8').Value = '''0.245'
$ws.Range('E8').Value = '  -0.92%  '

$ws.Range('E9').Value = '  -1.27%  '

$ws.Range('D10').Value = '''17.95'
$ws.Range('E10').Value = '  -1.65%  '

$ws.Range('E11').Value = '  +2.58%  '

$ws.Range('D12').Value = '1.816.55'
$ws.Range('E12').Value = '  +0.33%  '

$ws.Range('D13').Value = '1.612.91'
$ws.Range('E13').Value = '  +1.45%  '

$ws.Range('E14').Value = '  -1.07%  '

$ws.Range('E15').Value = '  -0.06%  '

$ws.Range('D16').Value = '25.984.42'
$ws.Range('E16').Value = '  +0.32%  '

$ws.Range('D17').Value = '''60.05'
$ws.Range('E17').Value = '  -0.30%  '

$ws.Range('D18').Value = '0.0₃0722'
$ws.Range('E18').Value = '  -0.34%  '

$ws.Range('E19').Value = '  -0.13%  '

$ws.Range('D20').Value = '''200.20'
$ws.Range('E20').Value = '  +3.51%  '

$ws.Range('D21').Value = '''4.22'
$ws.Range('E21').Value = '  +0.53%  '

$ws.Range('E22').Value = '  -1.95%  '

$ws.Range('D23').Value = '''6.00'
$ws.Range('E23').Value = '  +0.88%  '

$ws.Range('D24').Value = '''1.85'
$ws.Range('E24').Value = '  +8.48%  '

$ws.Range('D25').Value = '''143.14'
$ws.Range('E25').Value = '  +1.09%  '

$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('E27').Value = '  -8.31%  '

$ws.Range('E28').Value = '  -0.53%  '

$ws.Range('E29').Value = '  -0.54%  '

$ws.Range('E30').Value = '  +0.06%  '

$ws.Range('E31').Value = '  +0.24%  '

$ws.Range('E32').Value = '  -0.05%  '

$ws.Range('E33').Value = '  -3.13%  '

$ws.Range('E34').Value = '  -1.56%  '

$ws.Range('E35').Value = '  +0.10%  '

$ws.Range('D36').Value = '1.122.37'
$ws.Range('E36').Value = '  +1.18%  '

$ws.Range('E37').Value = '  +7.70%  '

$ws.Range('E38').Value = '  -0.07%  '

$ws.Range('E39').Value = '  -1.24%  '

$ws.Range('D40').Value = '''0.782'
$ws.Range('E40').Value = '  -0.28%  '

$ws.Range('D41').Value = '''0.490'
$ws.Range('E41').Value = '  -3.41%  '

$ws.Range('D42').Value = '''0.784'
$ws.Range('E42').Value = '  -4.47%  '

$ws.Range('D43').Value = '1.727.93'
$ws.Range('E43').Value = '  +0.22%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''5.09'
$ws.Range('E44').Value = '  -1.68%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '''92.50'
$ws.Range('E45').Value = '  -1.32%  '

$ws.Range('E46').Value = '  -1.09%  '

$ws.Range('D47').Value = '''53.31'
$ws.Range('E47').Value = '  -0.46%  '

$ws.Range('E48').Value = '  -1.38%  '

$ws.Range('E49').Value = '  -0.06%  '

$ws.Range('E50').Value = '  +0.12%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₇0914'
$ws.Range('E51').Value = '  -17.69%  '
